$wb = $excel.ActiveWorkbook

# CRED01_DealSetup
$ws = $wb.Worksheets.Item("CRED01_DealSetup")
$ws.Range("E2").Value = "UAT5_07092020114431"
$ws.Range("F2").Value = "UAT507092020114431"
$ws.Range("G2").Value = "FACILITY-A_07092020115134ZWO"
$ws.Range("G3").Value = "FACILITY-B_07092020115901PQX"
$ws.Range("G4").Value = "FACILITY-C_07092020120504VXK"

# CRED02_FacilitySetup
$ws = $wb.Worksheets.Item("CRED02_FacilitySetup")
$ws.Range("D2").Value = "UAT5_07092020114431"
$ws.Range("E2").Value = "FACILITY-A_07092020115134ZWO"
$ws.Range("D3").Value = "UAT5_07092020114431"
$ws.Range("E3").Value = "FACILITY-B_07092020115901PQX"
$ws.Range("D4").Value = "UAT5_07092020114431"
$ws.Range("E4").Value = "FACILITY-C_07092020120504VXK"

# CRED08_FacilityFeeSetup
$ws = $wb.Worksheets.Item("CRED08_FacilityFeeSetup")
$ws.Range("C2").Value = "FACILITY-A_07092020115134ZWO"
$ws.Range("C3").Value = "FACILITY-B_07092020115901PQX"
$ws.Range("C4").Value = "FACILITY-C_07092020120504VXK"

# SERV01_LoanDrawdown
$ws = $wb.Worksheets.Item("SERV01_LoanDrawdown")
$ws.Range("C2").Value = "UAT5_07092020114431"
$ws.Range("C3").Value = "UAT5_07092020114431"
$ws.Range("E3").Value = "FACILITY-A_07092020115134ZWO"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "60000667"
$ws.Range("C4").Value = "UAT5_07092020114431"
$ws.Range("E4").Value = "FACILITY-B_07092020115901PQX"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "60000668"
$ws.Range("C5").Value = "UAT5_07092020114431"

# SERV08C_ComprehensiveRepricing
$ws = $wb.Worksheets.Item("SERV08C_ComprehensiveRepricing")
$ws.Range("C2").Value = "UAT5_07092020114431"
$ws.Range("C3").Value = "UAT5_07092020114431"
$ws.Range("D3").Value = "FACILITY-A_07092020115134ZWO"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "60000667"
$ws.Range("C4").Value = "UAT5_07092020114431"
$ws.Range("D4").Value = "FACILITY-B_07092020115901PQX"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "60000668"
$ws.Range("C5").Value = "UAT5_07092020114431"
